$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 102, pushing the existing historical
# data (old rows 102-135) down to rows 104-137.
$ws.Rows.Item(102).Resize(2).Insert()

# Populate the two newly inserted rows (102-103) with the new weekly
# sample: same market/category/quality pairing as the row that used to
# be at 102-103, just with a new date.
$ws.Cells.Item(102, 1).Value = 11
$ws.Cells.Item(102, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(102, 3).Value = "Bíobío"
$ws.Cells.Item(102, 4).Value = 44663
$ws.Cells.Item(102, 5).Value = 8
$ws.Cells.Item(102, 6).Value = 100112044
$ws.Cells.Item(102, 7).Value = "Perejil"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 200
$ws.Cells.Item(102, 11).Value = 600
$ws.Cells.Item(102, 12).Value = 700
$ws.Cells.Item(102, 13).Value = 650
$ws.Cells.Item(102, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(102, 15).Value = "Región de Ñuble"
$ws.Cells.Item(102, 16).Value = 650
$ws.Cells.Item(102, 17).Value = 1
$ws.Cells.Item(102, 18).Value = "Hortaliza"

$ws.Cells.Item(103, 1).Value = 11
$ws.Cells.Item(103, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(103, 3).Value = "Bíobío"
$ws.Cells.Item(103, 4).Value = 44663
$ws.Cells.Item(103, 5).Value = 8
$ws.Cells.Item(103, 6).Value = 100112044
$ws.Cells.Item(103, 7).Value = "Perejil"
$ws.Cells.Item(103, 8).Value = "Sin especificar"
$ws.Cells.Item(103, 9).Value = "Segunda"
$ws.Cells.Item(103, 10).Value = 100
$ws.Cells.Item(103, 11).Value = 500
$ws.Cells.Item(103, 12).Value = 500
$ws.Cells.Item(103, 13).Value = 500
$ws.Cells.Item(103, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(103, 15).Value = "Región de Ñuble"
$ws.Cells.Item(103, 16).Value = 500
$ws.Cells.Item(103, 17).Value = 1
$ws.Cells.Item(103, 18).Value = "Hortaliza"
